$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "service": update unit counts / totals
# ----------------------------------------------------------------------
$service = $wb.Worksheets.Item("service")

$service.Range("B2").Value = 180
$service.Range("C2").Value = 9000
$service.Range("B3").Value = 44
$service.Range("C3").Value = 9000
$service.Range("B4").Value = 20
$service.Range("C4").Value = 9000
$service.Range("B5").Value = 180
$service.Range("C5").Value = 4000
$service.Range("B6").Value = 44
$service.Range("C6").Value = 4000
$service.Range("B7").Value = 20
$service.Range("C7").Value = 4000

# Move the selection on this sheet to B8 (was C5)
$service.Range("B8").Select() | Out-Null

# ----------------------------------------------------------------------
# Sheet "request": update the summary matrix
# ----------------------------------------------------------------------
$request = $wb.Worksheets.Item("request")

$request.Range("B2").Value = 500
$request.Range("C2").Value = 300
$request.Range("E2").Value = 500
$request.Range("F2").Value = 300

$request.Range("C3").Value = 300
$request.Range("D3").Value = 500
$request.Range("F3").Value = 300
$request.Range("G3").Value = 500

$request.Range("B4").Value = 300
$request.Range("C4").Value = 500
$request.Range("E4").Value = 300
$request.Range("F4").Value = 500

# ----------------------------------------------------------------------
# Sheet "intervalForSendingRequests": move its (inactive) selection
# ----------------------------------------------------------------------
$interval = $wb.Worksheets.Item("intervalForSendingRequests")
$interval.Range("D7").Select() | Out-Null

# ----------------------------------------------------------------------
# Make "request" the active sheet/tab, with a new selection
# ----------------------------------------------------------------------
$request.Activate() | Out-Null
$request.Range("E11").Select() | Out-Null
